$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.664.46'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.914.51'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Formula = "=""603.41"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Formula = "=""169.44"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("D7").Value = '3.918.58'
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Formula = "=""0.532"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").Formula = "=""6.45"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Formula = "=""0.0000254"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").Formula = "=""37.28"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").Value = '4.569.22'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '3.917.97'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '68.611.99'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Formula = "=""18.16"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +6.08%  '
$ws.Range("D19").Formula = "=""7.45"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").Formula = "=""10.89"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -1.13%  '
$ws.Range("D22").Formula = "=""473.25"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("D23").Formula = "=""0.743"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +2.76%  '
$ws.Range("D24").Formula = "=""0.0000168"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +0.83%  '
$ws.Range("D25").Formula = "=""83.90"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("E26").Value = '  +1.48%  '
$ws.Range("D27").Formula = "=""12.24"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +1.65%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Formula = "=""10.05"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").Formula = "=""1.00"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").Formula = "=""2.98"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").Value = '4.063.97'
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").Formula = "=""7.86"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +1.73%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Formula = "=""31.67"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Formula = "=""2.32"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("D35").Formula = "=""9.51"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  +2.68%  '
$ws.Range("D36").Value = '3.886.49'
$ws.Range("E36").Value = '  +0.80%  '
$ws.Range("E37").Value = '  -1.45%  '
$ws.Range("D38").Formula = "=""3.66"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +15.32%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").Formula = "=""1.04"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Formula = "=""0.141"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +2.52%  '
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("D42").Formula = "=""1.00"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("D44").Formula = "=""431.60"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Formula = "=""2.01"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +1.50%  '
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").Formula = "=""0.000299"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +13.08%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Formula = "=""1.00"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Formula = "=""8.66"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +1.70%  '
$ws.Range("D49").Formula = "=""47.28"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("D50").Formula = "=""143.77"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Formula = "=""26.41"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +0.88%  '
